# Added metadata for baseline sampling (Daily sheet rows for 2022-09-29 / 2022-09-30)

$wb = $excel.ActiveWorkbook

$daily = $wb.Worksheets.Item("Daily")
$weekly = $wb.Worksheets.Item("Weekly")

# --- New rows on the "Daily" sheet ---
# Row 2 : 2022-09-29
$daily.Range("A2").Value = 20220929
$daily.Range("B2").Value = "completed"
$daily.Range("C2").Value = "completed"
$daily.Range("D2").Value = "AH"
$daily.Range("E2").Value = "No mortality observed, tanks running well "

# Row 3 : 2022-09-30
$daily.Range("A3").Value = 20220930
$daily.Range("B3").Value = "completed "
$daily.Range("C3").Value = "completed"
$daily.Range("D3").Value = "AH"
$daily.Range("E3").Value = "Removed algae from a couple frags, no mortality observed "

# --- Selection / active-tab bookkeeping ---
# "Daily" becomes the active / selected sheet, with A4 selected next.
$daily.Select()
$daily.Range("A4").Select()

# "Weekly" is no longer the tab-selected sheet; its lingering selection moves to B17.
$weekly.Range("B17").Select()

# Re-activate "Daily" so it ends up as the workbook's active sheet.
$daily.Activate()
